# Remove the slide number from the template slide.
#
# Turning off the "Slide Number" header/footer element is the
# PowerPoint-native way to drop the slide-number placeholder from a
# slide: it deletes the <p:sp> for the sldNum placeholder outright
# (rather than merely clearing the <a:fld> inside it, which is what a
# plain Shape.Delete() on a placeholder does).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.HeadersFooters.SlideNumber.Visible = $false
